# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 90 (pushing the existing
# historical rows 90-133 down to 91-134) for the
# "Fruta, Agrícola del Norte S.A. de Arica - Naranja" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 90 - this shifts the
# previous rows 90..133 down to 91..134 (copying down their formatting,
# including the date number format used in column D).
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with this week's observation.
$ws.Range("A90").Value = 1
$ws.Range("B90").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C90").Value = "Arica y Parinacota"
$ws.Range("D90").Value = 45007
$ws.Range("E90").Value = 15
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100102
$ws.Range("H90").Value = "Cítricos"
$ws.Range("I90").Value = 100102005
$ws.Range("J90").Value = "Naranja"
$ws.Range("K90").Value = "Valencia"
$ws.Range("L90").Value = "Segunda"
$ws.Range("M90").Value = 300
$ws.Range("N90").Value = 1000
$ws.Range("O90").Value = 1100
$ws.Range("P90").Value = 1050
$ws.Range("Q90").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R90").Value = "Región de Coquimbo"
$ws.Range("S90").Value = 1050
$ws.Range("T90").Value = 1
